## CANumbers: fill in auto-filled Name/Address for the CA Number, drop the
## two other CA Numbers (kept blank), and do a test-run.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("CANumbers")
$ws1.Range("B2").Value = "Mr. RAMESHWAR DASS BOHET"
$ws1.Range("E2").Value = "PLOT NO 19 KH. NO. 38/22 GROUND FLOOR GRAM SABHA STREET NO. GALI NO-2 VILLAGE MAMURPUR CITY DELHI 110040 LANDMARK NEAR RADHA SWAMI ASHRAM"
$ws1.Range("A3").Clear()
$ws1.Range("A4").ClearContents()

$ws1.Rows.Item(1).RowHeight = 16.15
$ws1.Columns.Item(2).ColumnWidth = 28.15
$ws1.Columns.Item(3).ColumnWidth = 10.1
$ws1.Columns.Item(4).ColumnWidth = 8.1
$ws1.Columns.Item(5).ColumnWidth = 145.95
$ws1.Range("B17").Select() | Out-Null

## ComplaintNumbers: log two test complaints against the CA Number.
$ws2 = $wb.Worksheets.Item("ComplaintNumbers")
$ws2.Columns.Item(1).ColumnWidth = 16.166
$ws2.Columns.Item(2).ColumnWidth = 13.5
$ws2.Columns.Item(3).ColumnWidth = 21.5

$ws2.Range("A2").NumberFormat = "YYYY\-MM\-DD\ H:MM:SS"
$ws2.Range("A2").Value = 43677.0261515908
$ws2.Rows.Item(2).RowHeight = 12.8
$ws2.Range("B2").Value = 60014463164
$ws2.Range("C2").Value = "some complaint number"

$ws2.Range("A3").NumberFormat = "YYYY\-MM\-DD\ H:MM:SS"
$ws2.Range("A3").Value = 43677.0271245861
$ws2.Range("B3").Value = 60014463164
$ws2.Range("C3").Value = "some complaint number"

## Remarks: rewrite the canned remark phrases (one less row than before).
$ws3 = $wb.Worksheets.Item("Remarks")
$ws3.Range("A1").Value = "No supply in area."
$ws3.Range("A2").Value = "No electricity."
$ws3.Range("A3").Value = "Power breakdown in area, get it fixed immedieately."
$ws3.Rows.Item(4).Delete()
$ws3.Columns.Item(1).ColumnWidth = 42.166
$ws3.Rows.Item(1).RowHeight = 12.8
$ws3.Rows.Item(2).RowHeight = 12.8
$ws3.Rows.Item(3).RowHeight = 12.8
$ws3.PageSetup.CenterHeader = "&""Times New Roman,Regular""&12 &A"
$ws3.PageSetup.CenterFooter = "&""Times New Roman,Regular""&12 Page &P"
$ws3.Range("A5").Select() | Out-Null

## Sheet4: same header/footer touch-up as Remarks.
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.PageSetup.CenterHeader = "&""Times New Roman,Regular""&12 &A"
$ws4.PageSetup.CenterFooter = "&""Times New Roman,Regular""&12 Page &P"

## Switch the active tab to ComplaintNumbers and leave a selection there,
## matching the "testrun" the commit message refers to.
$ws2.Range("C15").Select() | Out-Null
$ws2.Activate() | Out-Null
